$d = $word.ActiveDocument

# The title paragraph currently holds several runs that together spell out
# "Relação de Alunos do {turma}  – 25/26". Find that paragraph (the one
# that contains the {turma} placeholder) so we target it regardless of
# exact paragraph index.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "*{turma}*") {
        $target = $candidate
        break
    }
}
if ($target -eq $null) {
    $target = $d.Paragraphs.Item($d.Paragraphs.Count)
}

# Replace the whole paragraph's contents with a single run that merges all
# the previous text fragments into one string and adds the {professor}
# placeholder, while switching the font references to the theme major font
# (majorHAnsi) and explicitly turning bold off, per the template update.
$paraXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="535B5CD2" w14:textId="337464F6" w:rsidR="00BF2514" w:rsidRPr="00903062" w:rsidRDefault="006E036C" w:rsidP="00903062"><w:pPr><w:pStyle w:val="TITULO"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:eastAsia="Abadi Extra Light" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:b w:val="0"/><w:bCs w:val="0"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r w:rsidRPr="004C50F1"><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:eastAsia="Abadi Extra Light" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:b w:val="0"/><w:bCs w:val="0"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>Relação de Alunos do {turma} - {professor} – 25/26</w:t></w:r></w:p>
'@

$target.Range.InsertXML($paraXml)
